$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row data: the long "Montaje de tubería..." boilerplate sentence is
# dropped from the description (col A), the placeholder "unidad suerte m"
# in col B is replaced by the real measured quantity for that pipe size,
# and the page reference in col D is updated. Col C ("m") is unchanged.
$rows = @(
    @{ Row = 2; Desc = "Montaje de tubería PE100, ø 25 mm, 1,0 MPa, colocada"; Qty = "109.449" },
    @{ Row = 3; Desc = "Montaje de tubería PE100, ø 32 mm, 1,0 MPa, colocada"; Qty = "92.249" },
    @{ Row = 4; Desc = "Montaje de tubería PE100, ø 40 mm, 1,0 MPa, colocada"; Qty = "67.386" },
    @{ Row = 5; Desc = "Montaje de tubería PE100, ø 50 mm, 1,0 MPa, colocada"; Qty = "67.383" },
    @{ Row = 6; Desc = "Montaje de tubería PE100, ø 63 mm, 1,0 MPa, colocada"; Qty = "32.687" },
    @{ Row = 7; Desc = "Montaje de tubería PE100, ø 75 mm, 1,0 MPa, colocada"; Qty = "12.050" },
    @{ Row = 8; Desc = "Montaje de tubería PE100, ø 90 mm, 1,0 MPa, colocada"; Qty = "7.956" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column A - shortened description
    $ws.Range("A$rowNum").Value = $r.Desc

    # Column B - quantity, kept as literal text (e.g. "12.050" would lose its
    # trailing zero if stored as a real number), so force the "@" text
    # format before assigning, then drop the leftover formatting again so
    # the cell is left with no explicit style, same as before the edit.
    $cellB = $ws.Range("B$rowNum")
    $cellB.NumberFormat = "@"
    $cellB.Value = $r.Qty
    $cellB.ClearFormats()

    # Column D - updated page/section reference
    $ws.Range("D$rowNum").Value = "Pág. 3, apartado 2"

    # Column E (referencia) - this column is no longer populated for data rows
    $ws.Range("E$rowNum").ClearContents()
}
